# Restore C10 value from 18 to 1, as per the admin revision restore.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
